# Add a new worksheet "ValidLogIn" with username/password sample data,
# then make it the active sheet (matching the commit "Updated Loginpage with errMsg").

$wb = $excel.ActiveWorkbook
$tc1 = $wb.Worksheets.Item(1)

# Add the new sheet right after the existing "TC1" sheet.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tc1)
$newSheet.Name = "ValidLogIn"

# Fill in the header row and the sample credentials row.
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "pointofsale"

# Column widths to match the authored layout (closest values the host's
# pixel-quantized ColumnWidth rounding can reach).
$newSheet.Columns.Item(1).ColumnWidth = 10.083
$newSheet.Columns.Item(2).ColumnWidth = 16.08

# Match the original selection (A2) carried over from the TC1 sheet.
$newSheet.Range("A2").Select()

# Make the new sheet the active / visible tab.
$newSheet.Activate()
